$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# The document currently looks like:
#   1. "Aktuelles"                         (heading)
#   2. (empty paragraph)
#   3. "Wir suchen ... in Vollzeit"         (bold job heading)
#   4. (empty, bordered paragraph)
#   5. "Haben Sie Interesse ..."
#   6-11. bulleted list items
#   12. "Dann schicken Sie uns eine Bewerbung zu!"
#   13. (empty, bordered paragraph)
#   14. "Wir bieten Ihnen eine kreative, ..."
#   15. "Waldkindergarten Mäusewiese Plankenfels / Scherleithen 5 / 95515 Plankenfels"
#   16. "Frau Annalena Hofmann, ... / Telefon ... / Mail ..."
#
# The whole job-ad block (paragraphs 2 through 16) is replaced by a single
# new paragraph with a short "news" text, containing two manual line breaks.
# --------------------------------------------------------------------------

# Remove paragraphs 2..16 completely (including paragraph 2's own mark), so
# the document body ends right after paragraph 1 ("Aktuelles").
$delStart = $d.Paragraphs(2).Range.Start
$delEnd = $d.Paragraphs($d.Paragraphs.Count).Range.End
$d.Range($delStart, $delEnd).Delete()

# New paragraph text (German copy about winter / Fasching).
$para1Text = "Der Winter und die kalte Jahreszeit haben begonnen und somit wollen wir gemeinsam Kraft für das kommende Jahr sammeln. Wir erkunden unsere Lieblingsplätze, gehen auf Entdeckungsreisen und machen oft Feuer, um uns daran zu wärmen. Wir wünschen uns viel Schnee, damit wir endlich zu dem Schlittenhang gehen können. "
$para2Text = "Im Februar ist unser Faschingsmonat, wir wollen uns bis zum Aschermittwoch verkleiden, denn das macht uns so viel Spaß! Es gibt viele Spiele, Staffelläufe, Musik und bunte Wimpel auf der Waldwiese."

# Build the replacement paragraph as raw WordprocessingML so the two
# <w:br/> line breaks land in their own runs, matching a hand-edited
# document rather than a single merged run.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">' + $para1Text + '</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>' + $para2Text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Paragraphs(1).Range.End
$d.Range($insertionPoint, $insertionPoint).InsertXML($xml)
